$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters in the order values are stored per row: B,C,D,E,F,I,J,K,L,M,N
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$values = @{}
$values[2] = @(1.02,1.024277781602758,1.034723261155134,0.9926147277508489,1.043143036014878,1.033020092456124,1.02945329371916,1.037521916480203,0.9955398523336033,1.045917733870858,1.013865946968037)
$values[3] = @(1.02,1.0250911147939,1.035345828564626,0.9936372048519304,1.043897636448387,1.033133606290964,1.02990616939188,1.037954227958405,0.9963617723202692,1.046483456819377,1.014016713643455)
$values[4] = @(1.02,1.025618015520876,1.035749046332516,0.9942998659930995,1.044386633226335,1.033205885101957,1.030199175446401,1.038233645832237,0.9968940712668345,1.046849572375537,1.014114227644884)
$values[5] = @(1.02,1.025839671246092,1.035918646949525,0.9945786998346017,1.044592377715559,1.033235989768147,1.030322345567013,1.038351035893999,0.997117960005301,1.047003498804661,1.014155211982999)
$values[6] = @(1.02,1.025876896768096,1.035947128736262,0.9946255319796338,1.044626933043065,1.033241027955128,1.030343025770135,1.038370741632611,0.9971555583673453,1.047029344350973,1.01416209280164)
$values[7] = @(1.02,1.025620976721397,1.035751312201414,0.9943035907982488,1.044389381727955,1.033206288468367,1.030200821291028,1.038235214709175,0.9968970624462087,1.046851629105184,1.014114775321646)
$values[8] = @(1.02,1.024552522004941,1.034933581924727,0.9929600610674301,1.043397906462511,1.033058697278089,1.029606351597084,1.037668082856387,0.995817528259106,1.046108910176356,1.013916907713642)
$values[9] = @(1.02,1.022674590776461,1.033495593468632,0.9906006454969559,1.041656405687771,1.032789676328352,1.028558610878733,1.03666636264437,0.9939188001724441,1.044800641281617,1.013567937647762)
$values[10] = @(1.02,1.021425981102837,1.032539036622451,0.989033133672735,1.040499298393637,1.032604356429939,1.027860053676367,1.035997046557271,0.9926553831429383,1.043928898153579,1.013335115525263)
$values[11] = @(1.02,1.020886134088841,1.032125359690579,0.988355674866747,1.03999920791617,1.032522703198482,1.02755757195652,1.035706886867343,0.9921088820399291,1.043551548783923,1.013234265210222)
$values[12] = @(1.02,1.020685734174174,1.031971781515322,0.9881042295826724,1.039813596101326,1.03249216264895,1.027445217812993,1.035599058792219,0.9919059725120875,1.043411404372365,1.013196799854351)
$values[13] = @(1.02,1.020728715028663,1.032004720910188,0.9881581567098651,1.039853403912612,1.032498723232991,1.027469318090278,1.035622190512851,0.9919494934313052,1.043441464893821,1.01320483652173)
$values[14] = @(1.02,1.020869566444035,1.03211266323035,0.9883348863814464,1.039983862234626,1.032520183005416,1.027548284699266,1.035697974773524,0.9920921077337197,1.043539963986183,1.013231168412693)
$values[15] = @(1.02,1.020956366040293,1.032179180696573,0.9884438009545853,1.040064261065044,1.032533377153811,1.027596938809689,1.035744661395731,0.9921799884222134,1.043600655222987,1.013247391703529)
$values[16] = @(1.02,1.021461826167857,1.032566502089342,0.9890781214508737,1.040532507836093,1.032609745870953,1.027880128465123,1.036016296430655,0.9926916645766087,1.043953944289978,1.013341807893359)
$values[17] = @(1.02,1.021779105820453,1.032809598765559,0.989476357848556,1.040826481074927,1.03265727335459,1.028057766287642,1.036186595528224,0.9930127773699352,1.044175586856561,1.013401023154667)
$values[18] = @(1.02,1.021964247635334,1.032951442840052,0.9897087662937556,1.040998041690161,1.032684859385194,1.028161379158533,1.036285895072883,0.9932001317071769,1.044304878788895,1.013435558832622)
$values[19] = @(1.02,1.0220273893808,1.03299981640749,0.9897880325774034,1.041056554796324,1.032694242428378,1.028196708384935,1.036319748003704,0.9932640239640975,1.044348965906557,1.013447333984417)
$values[20] = @(1.02,1.021745056640046,1.032783511617724,0.9894336180360679,1.040794931100697,1.032652188161881,1.028038707448437,1.036168327453562,0.9929783193494215,1.044151805497126,1.013394670284658)
$values[21] = @(1.02,1.020828085795185,1.032080874687326,0.9882828385668249,1.039945441516611,1.032513869453609,1.027525030962003,1.035675659562046,0.9920501090198102,1.043510957893917,1.013223414462976)
$values[22] = @(1.02,1.020252264080576,1.031639561951374,0.9875604150241495,1.039412168212454,1.032425683117767,1.027202068805302,1.035365612020102,0.9914670000341481,1.043108147832231,1.013115710096058)
$values[23] = @(1.02,1.020557449786283,1.031873465643722,0.9879432794643023,1.039694786665835,1.032472547756922,1.027373276088232,1.035530000875299,0.991776070289318,1.04332167346558,1.013172808814125)
$values[24] = @(1.02,1.021760441758365,1.032795299123006,0.9894529299347244,1.040809186895105,1.032654486361716,1.028047319318482,1.036176582112425,0.9929938892766442,1.044162551234261,1.013397540884287)
$values[25] = @(1.02,1.023159499008532,1.033866985783518,0.9912096547607049,1.042105948281079,1.032860280514781,1.028829493792954,1.036925602398836,0.9944092447426414,1.044800641281617,1.013658187715893)

foreach ($r in $values.Keys) {
    $rowVals = $values[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}
